$wb = $excel.ActiveWorkbook

# --- The previously-active sheet ("no_feot") loses its tab selection and
#     ends up with the whole used range selected once focus moves away. ---
$wsPrev = $wb.Worksheets.Item("no_feot")
[void]$wsPrev.Range("A1:XFD7").Select()

# --- Add the new "Feo_Fe2O3" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Feo_Fe2O3"
$ws.Activate()

# --- Header row (row 1) ---
# Same layout as "no_feot" but with an extra Fe2O3_Liq column inserted
# right after FeO_Liq (so columns shift one place to the right from F on).
$ws.Cells.Item(1,1).Value  = "Sample_ID"
$ws.Cells.Item(1,2).Value  = "SiO2_Liq"
$ws.Cells.Item(1,3).Value  = "TiO2_Liq"
$ws.Cells.Item(1,4).Value  = "Al2O3_Liq"
$ws.Cells.Item(1,5).Value  = "FeO_Liq"
$ws.Cells.Item(1,6).Value  = "Fe2O3_Liq"
$ws.Cells.Item(1,7).Value  = "MnO_Liq"
$ws.Cells.Item(1,8).Value  = "MgO_Liq"
$ws.Cells.Item(1,9).Value  = "CaO_Liq"
$ws.Cells.Item(1,10).Value = "Na2O_Liq"
$ws.Cells.Item(1,11).Value = "K2O_Liq"
$ws.Cells.Item(1,12).Value = "Cr2O3_Liq"
$ws.Cells.Item(1,13).Value = "P2O5_Liq"
$ws.Cells.Item(1,14).Value = "H2O_Liq"

# --- Data rows 2-6 ---
$sampleName = "Sisson, T.W., Grove, T.L. (1993)"

$data = @(
    # B,    C,    D,    E,    F,    G,    H,    I,     J,    K,    L, M,    N
    @(51.1, 0.93, 17.5, 8.91, 1,    0.18, 6.09, 11.5,  3.53, 0.17, 0, 0.15, 3.8),
    @(51.5, 1.19, 19.2, 8.6999999999999993, 1.5, 0.19, 4.9800000000000004, 10, 3.72, 0.42, 0, 0.14000000000000001, 6.2),
    @(59.1, 0.54, 19.100000000000001, 5.22, 2.2000000000000002, 0.19, 3.25, 7.45, 4, 0.88, 0, 0.31, 6.2),
    @(52.5, 0.98, 19.2, 8.0399999999999991, 1.2, 0.2, 4.99, 9.64, 4.1500000000000004, 0.21, 0, 0.14000000000000001, 6.2),
    @(56.2, 0.34, 20.399999999999999, 5.88, 1.3, 0.2, 2.58, 7.18, 6.02, 1.02, 0, 0.23, 6.2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row,1).Value = $sampleName
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row,$col).Value = $vals[$j]
    }
}

# --- Selection / active sheet bookkeeping to mirror the authored workbook ---
[void]$ws.Range("F7").Select()
$ws.Activate()

Write-Host "done"
